# DailyLog.xlsx — "8/25 morning at home"
# Append the new LeetCode log entry for 8/25/2021 as row 34:
#   Day = 8/25/2021 (serial 44433), Question Number = 348_DesignTicTacToe, Topic = design

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(34, 1).Value = 44433
$ws.Cells.Item(34, 1).NumberFormat = "d-mmm"
$ws.Cells.Item(34, 2).Value = "348_DesignTicTacToe"
$ws.Cells.Item(34, 3).Value = "design"

# Leave the selection where the author's session ended up.
$ws.Range("I24").Select() | Out-Null
